# "D suite.xlsx" - update the Runmode column on the "Test Cases" sheet.
#
# The diff shows every "N" runmode value (shared-string driven) in column C
# (rows 2-5) turning into "Y" - i.e. these test cases were switched from
# "don't run" to "run". The "Results" column (D) keeps showing "SKIP" (its
# shared-string index just shifts because the string table gets re-ordered
# by the save), and the selection on the sheet ends up on C6 (just past the
# edited range) instead of the old C2:C5 selection.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test Cases")

$ws1.Range("C2:C5").Value = "Y"

# Leave the selection where the author's session ended up - on C6, right
# below the block that was just edited.
$ws1.Range("C6").Select()
